# Add sanunits (Biogenic Refinery) impact items to the "info" and "GWP"
# sheets of the workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("info")
$ws2 = $wb.Worksheets.Item("GWP")

# ---------------------------------------------------------------------
# "info" sheet: ID (col A) / functional_unit (col B), rows 14-23
# ---------------------------------------------------------------------
$infoRows = @(
    @(14, "ElectricMotor",          "ea"),
    @(15, "Electronics",            "kg"),
    @(16, "CatalyticConverter",     "ea"),
    @(17, "OilHeatExchanger",       "ea"),
    @(18, "Pump",                   "ea"),
    @(19, "HydronicHeatExchanger",  "ea"),
    @(20, "ElectricConnectors",     "kg"),
    @(21, "ElectricCables",         "m"),
    @(22, "PVC",                    $null),
    @(23, "PE",                     $null)
)

# Row 14 introduces both a new item name and a new functional-unit string;
# the functional unit ("ea") was registered first in the shared-string
# table, so write column B before column A for that row to reproduce the
# original authoring order. For the remaining rows, column A is written
# first (matching natural top-to-bottom, left-to-right entry).
$ws1.Range("B14").Value = "ea"
$ws1.Range("A14").Value = "ElectricMotor"

foreach ($row in $infoRows) {
    $r = $row[0]
    if ($r -eq 14) { continue }
    $ws1.Range("A$r").Value = $row[1]
    if ($row[2] -ne $null) {
        $ws1.Range("B$r").Value = $row[2]
    }
}

# ---------------------------------------------------------------------
# "GWP" sheet: ID, unit, expected, low, high, distribution, references
# rows 14-24 (note row 24 "Polyacrylamide" has no counterpart on "info")
# ---------------------------------------------------------------------
$gwpRows = @(
    @(14, "ElectricMotor",         9.9703471209999996,    8.9733124089,          10.967381833100001),
    @(15, "Electronics",           8.8474600330000008,    7.9627140297000008,    9.7322060363000009),
    @(16, "CatalyticConverter",    13.937562,              12.543805799999999,    15.3313182),
    @(17, "OilHeatExchanger",      114014.66383541599,     102613.1974518744,     125416.13021895761),
    @(18, "Pump",                  8.5937041260028693,     7.7343337134025827,    9.4530745386031576),
    @(19, "HydronicHeatExchanger", 324.96867299992499,     292.47180569993247,    357.46554029991751),
    @(20, "ElectricConnectors",    8.8474600327818393,     7.9627140295036556,    9.732206036060024),
    @(21, "ElectricCables",        4.3790218750000003,     3.9411196875000005,    4.8169240625000009),
    @(22, "PVC",                   1,                       0,                     2),
    @(23, "PE",                    1,                       0,                     2),
    @(24, "Polyacrylamide",        2.797495751,             2.5177461759000002,    3.0772453261000003)
)

foreach ($row in $gwpRows) {
    $r = $row[0]

    $ws2.Range("A$r").Value = $row[1]
    $ws2.Range("B$r").Value = "kg CO2-eq"
    $ws2.Range("C$r").Value = $row[2]
    $ws2.Range("D$r").Value = $row[3]
    $ws2.Range("E$r").Value = $row[4]
    $ws2.Range("F$r").Value = "uniform"
    $ws2.Range("G$r").Value = "ecoinvent 3"

    # Match the shaded-fill formatting used by the other "uniform"
    # rows (e.g. row 5, "Concrete") for columns C:F.
    $ws2.Range("C5:F5").Copy()
    $ws2.Range("C$r").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Restore the view/selection state recorded in the saved workbook.
# ---------------------------------------------------------------------
$ws1.Range("A22:A23").Select()
$ws2.Activate()
$ws2.Range("E28").Select()
